# Generate Report for Handoff
# A new handoff was generated for file "2ad9d8b4-4841-443e-9fbc-ef9bbdeee56e.md",
# refreshing its "Latest Handoff" timestamp on every sheet of the report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G5").Value = "2016-08-19 06:40:51"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H5").Value = "2016-08-19 06:40:46"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H5").Value = "2016-08-19 06:40:51"
